# Helper: Excel COM colors are BGR-packed longs (same as the VBA RGB() function).
function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing sheet and add the new "Instalación" sheet after it.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Presupuesto"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Instalación"

# ---------------------------------------------------------------------------
# 2. Build the "Instalación" sheet content.
# ---------------------------------------------------------------------------

# -- Title bar (merged C6:G7) --
$title = $ws2.Range("C6:G7")
$title.Merge()
$title.Value = "Tabla de instalación y funcionamiento"
$title.Font.Name = "Calibri"
$title.Font.Bold = $true
$title.Font.Size = 14
$title.Font.Color = RGB 255 255 255
$title.Interior.Color = RGB 0x49 0x50 0x57
$title.HorizontalAlignment = -4131
$title.VerticalAlignment = -4108
$title.Borders.LineStyle = 1
$title.Borders.Weight = -4138
$title.Borders.Color = RGB 0x20 0x26 0x16

$ws2.Rows.Item(5).RowHeight = 15
$ws2.Rows.Item(6).RowHeight = 15
$ws2.Rows.Item(7).RowHeight = 15

# -- Header row (row 8) --
$headers = @("Fecha", "Conexión física", "Conexión local", "Conexión Internet", "Observaciones")
$headerRange = $ws2.Range("C8:G8")
for ($i = 0; $i -lt 5; $i++) {
    $cell = $ws2.Cells.Item(8, 3 + $i)
    $cell.Value = $headers[$i]
}
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 12
$headerRange.Font.Color = RGB 0 0 0
$headerRange.Interior.Color = RGB 0x59 0x83 0x92
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = -4138
$headerRange.Borders.Color = RGB 0x20 0x26 0x16
$ws2.Rows.Item(8).RowHeight = 16.2

# -- Data rows 9, 10, 11, 12 --
$dates = 45509, 45510, 45511
for ($i = 0; $i -lt 3; $i++) {
    $row = 9 + $i
    $c = $ws2.Cells.Item($row, 3)
    $c.Value = $dates[$i]
    $c.NumberFormat = "mm-dd-yy"
}

$ws2.Range("D9").Value = "✔"
$ws2.Range("E9").Value = "✘"
$ws2.Range("F9").Value = "✘"
$ws2.Range("G9").Value = "El primer día de instalación se verificó el encendido con una red personal del autor pero no se logró conectar a la red del Instituto."

$ws2.Range("D10").Value = "✔"
$ws2.Range("E10").Value = "✘"
$ws2.Range("F10").Value = "✘"
$ws2.Range("G10").Value = "El segundo día se hicieron pruebas de conexión local sin éxito y la colocación del módulo prototipo en el puesto 3 del parqueadero."

$ws2.Range("D11").Value = "✔"
$ws2.Range("E11").Value = "✔"
$ws2.Range("F11").Value = "✘"
$ws2.Range("G11").Value = "Se instaló el módulo en el puesto 2 del parqueadero y fuente de energía para alimentación y se estableció conexión con la red del Instituto, la conexión suele ser debil y presenta algunos retrasos de la información en la página."

# Row 12 left blank intentionally (trailing styled row in the source table).

# -- Formatting for the C (date) column --
$cDataFull = $ws2.Range("C10:C12")
$cDataFull.Font.Name = "Calibri"
$cDataFull.Font.Bold = $true
$cDataFull.Font.Color = RGB 255 255 255
$cDataFull.Interior.Color = RGB 0x12 0x45 0x59
$cDataFull.HorizontalAlignment = -4108
$cDataFull.VerticalAlignment = -4108
$cDataFull.Borders.LineStyle = 1
$cDataFull.Borders.Weight = 2
$cDataFull.Borders.Color = RGB 0x20 0x26 0x16

$cFirst = $ws2.Range("C9")
$cFirst.Font.Name = "Calibri"
$cFirst.Font.Bold = $true
$cFirst.Font.Color = RGB 255 255 255
$cFirst.Interior.Color = RGB 0x12 0x45 0x59
$cFirst.HorizontalAlignment = -4108
$cFirst.VerticalAlignment = -4108
$cFirst.Borders.Item(7).LineStyle = 1
$cFirst.Borders.Item(7).Weight = 2
$cFirst.Borders.Item(7).Color = RGB 0x20 0x26 0x16
$cFirst.Borders.Item(10).LineStyle = 1
$cFirst.Borders.Item(10).Weight = 2
$cFirst.Borders.Item(10).Color = RGB 0x20 0x26 0x16
$cFirst.Borders.Item(9).LineStyle = 1
$cFirst.Borders.Item(9).Weight = 2
$cFirst.Borders.Item(9).Color = RGB 0x20 0x26 0x16

# -- D:F check/cross columns (rows 9-12) --
$checks = $ws2.Range("D9:F12")
$checks.Interior.Color = RGB 0xCC 0xCC 0xCC
$checks.HorizontalAlignment = -4108
$checks.VerticalAlignment = -4108
$checks.Borders.LineStyle = 1
$checks.Borders.Weight = 2
$checks.Borders.Color = RGB 0x20 0x26 0x16

# -- Observations column G (rows 9-12) --
$obs = $ws2.Range("G9:G12")
$obs.Interior.Color = RGB 0xCC 0xCC 0xCC
$obs.WrapText = $true
$obs.VerticalAlignment = -4108
$obs.Borders.LineStyle = 1
$obs.Borders.Weight = 2
$obs.Borders.Color = RGB 0x20 0x26 0x16

$ws2.Rows.Item(9).RowHeight = 57.6
$ws2.Rows.Item(10).RowHeight = 57.6
$ws2.Rows.Item(11).RowHeight = 100.8

# -- Column widths --
$ws2.Columns.Item(4).ColumnWidth = 15
$ws2.Columns.Item(5).ColumnWidth = 14.67
$ws2.Columns.Item(6).ColumnWidth = 18
$ws2.Columns.Item(7).ColumnWidth = 30.1

$ws2.Range("A1").Select()
$ws2.Range("I10").Select()

# ---------------------------------------------------------------------------
# 3. Update the view on "Presupuesto": clear old selection, zoom, scroll.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A26").Select()
$excel.ActiveWindow.Zoom = 115

# ---------------------------------------------------------------------------
# 4. Activate "Instalación" last so it becomes the active/visible tab.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("M7").Select()
